$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Dataset" sheet had TestCase07/08/09/10 rows jumbled (a stray
# "TestCase09" label on rows that belong to TestCase08, a duplicated blank
# "[]" row mislabeled "TestCase10", etc). This change re-sorts/relabels rows
# 20-31 so that:
#   TestCase07 -> 3 parameter rows
#   TestCase08 -> 3 parameter rows
#   TestCase11 -> 3 parameter rows (shifted up, was rows 28-30)
#   TestCase12 -> 3 parameter rows (shifted up, was rows 31-33)
# and the two now-unused trailing rows (32-33) are cleared out entirely.
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 20; A = "TestCase07"; B = "[Tan Leatherette Weekender Duffle;Xtend Smart Watch]" }
    @{ Row = 21; A = "TestCase07"; B = "[Connector;Xtend Smart Watch]" }
    @{ Row = 22; A = "TestCase07"; B = "[Kindle;Jenga]" }
    @{ Row = 23; A = "TestCase08"; B = "[Tan Leatherette Weekender Duffle, 60.0]" }
    @{ Row = 24; A = "TestCase08"; B = "[SuitCase, 50.0]" }
    @{ Row = 25; A = "TestCase08"; B = "[Jenga, 60.0]" }
    @{ Row = 26; A = "TestCase11"; B = "[crio user, criouser@gmail.com, Testing the contact us page]" }
    @{ Row = 27; A = "TestCase11"; B = "[facebook user, test_user__@gmail.com, !!!special characters!!]" }
    @{ Row = 28; A = "TestCase11"; B = "[hacker user !!!, bad_user@gmail.com, <XSS testing>]" }
    @{ Row = 29; A = "TestCase12"; B = "[Yarine Floor Lamp, Addr line 1 addr Line 2 addr line 3]" }
    @{ Row = 30; A = "TestCase12"; B = "[Connector, TEST ADDR LINES COUNT GREATER TH]" }
    @{ Row = 31; A = "TestCase12"; B = "[Connector, 1 Hacker Way Menlo Park, CA 94025]" }
)

# Rows in column A that need to end up on the "style 2" cell format
# (matches the formatting already used by A29:A33 in the original sheet).
$styleTwoRows = @(20, 22, 23, 24, 25, 27, 28, 29, 30, 31)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}

# Copy the pre-existing "style 2" formatting onto the column-A cells that
# need it, using a still-valid style-2 source cell (A29 keeps style 2 both
# before and after the value rewrite above).
$styleSource = $ws.Range("A29")
$styleSource.Copy()
foreach ($r in $styleTwoRows) {
    $ws.Range("A" + $r).PasteSpecial(-4122)
}

# The former rows 32/33 (old TestCase12 rows that are now rows 29-31) are no
# longer part of the table - blank them out completely (formatting removed
# too, same as the untouched rows below them).
$ws.Range("A32:B33").Clear()

# Two trailing blank formatted rows at the very end of the sheet are dropped.
$ws.Rows("998:999").Delete()
